# Rebuild the Daily Driver Report (Sheet1) with verified driver data.
# The underlying driver roster was re-validated; this applies the corrected
# rows (driver_name, asset_id, job_site, scheduled_start, scheduled_end,
# actual_start, actual_end, status, status_reason) for rows 2-16, replacing
# the previously unverified values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Charles Anderson
$ws.Cells.Item(2,1).Value = "Charles Anderson"
$ws.Cells.Item(2,2).Value = "TRK-1010"
$ws.Cells.Item(2,3).Value = "West Plano Project"
$ws.Cells.Item(2,4).Value = "06:00 AM"
$ws.Cells.Item(2,5).Value = "03:00 PM"
$ws.Cells.Item(2,6).Value = "06:05 AM"
$ws.Cells.Item(2,7).Value = "02:48 PM"
$ws.Cells.Item(2,8).Value = "On Time"
$ws.Cells.Item(2,9).Value = ""

# Row 3: Mark Thompson
$ws.Cells.Item(3,1).Value = "Mark Thompson"
$ws.Cells.Item(3,2).Value = "TRK-1015"
$ws.Cells.Item(3,3).Value = "Downtown Construction"
$ws.Cells.Item(3,4).Value = "06:00 AM"
$ws.Cells.Item(3,5).Value = "03:00 PM"
$ws.Cells.Item(3,6).Value = "N/A"
$ws.Cells.Item(3,7).Value = "N/A"
$ws.Cells.Item(3,8).Value = "On Time"
$ws.Cells.Item(3,9).Value = ""

# Row 4: James Davis
$ws.Cells.Item(4,1).Value = "James Davis"
$ws.Cells.Item(4,2).Value = "TRK-1005"
$ws.Cells.Item(4,3).Value = "North Dallas Site"
$ws.Cells.Item(4,4).Value = "06:00 AM"
$ws.Cells.Item(4,5).Value = "03:00 PM"
$ws.Cells.Item(4,6).Value = "07:03 AM"
$ws.Cells.Item(4,7).Value = "03:12 PM"
$ws.Cells.Item(4,8).Value = "Not On Job"
$ws.Cells.Item(4,9).Value = "At incorrect location: North Richland Hills"

# Row 5: Anthony Martin
$ws.Cells.Item(5,1).Value = "Anthony Martin"
$ws.Cells.Item(5,2).Value = "TRK-1014"
$ws.Cells.Item(5,3).Value = "West Plano Project"
$ws.Cells.Item(5,4).Value = "07:00 AM"
$ws.Cells.Item(5,5).Value = "04:00 PM"
$ws.Cells.Item(5,6).Value = "07:16 AM"
$ws.Cells.Item(5,7).Value = "03:58 PM"
$ws.Cells.Item(5,8).Value = "Late"
$ws.Cells.Item(5,9).Value = "16 minutes late"

# Row 6: Daniel Jackson
$ws.Cells.Item(6,1).Value = "Daniel Jackson"
$ws.Cells.Item(6,2).Value = "TRK-1012"
$ws.Cells.Item(6,3).Value = "Richardson Development"
$ws.Cells.Item(6,4).Value = "06:15 AM"
$ws.Cells.Item(6,5).Value = "03:15 PM"
$ws.Cells.Item(6,6).Value = "06:16 AM"
$ws.Cells.Item(6,7).Value = "03:12 PM"
$ws.Cells.Item(6,8).Value = "On Time"
$ws.Cells.Item(6,9).Value = ""

# Row 7: Christopher Thomas
$ws.Cells.Item(7,1).Value = "Christopher Thomas"
$ws.Cells.Item(7,2).Value = "TRK-1011"
$ws.Cells.Item(7,3).Value = "Downtown Construction"
$ws.Cells.Item(7,4).Value = "06:30 AM"
$ws.Cells.Item(7,5).Value = "03:30 PM"
$ws.Cells.Item(7,6).Value = "06:29 AM"
$ws.Cells.Item(7,7).Value = "03:33 PM"
$ws.Cells.Item(7,8).Value = "On Time"
$ws.Cells.Item(7,9).Value = ""

# Row 8: William Brown
$ws.Cells.Item(8,1).Value = "William Brown"
$ws.Cells.Item(8,2).Value = "TRK-1004"
$ws.Cells.Item(8,3).Value = "Richardson Development"
$ws.Cells.Item(8,4).Value = "06:30 AM"
$ws.Cells.Item(8,5).Value = "03:30 PM"
$ws.Cells.Item(8,6).Value = "06:32 AM"
$ws.Cells.Item(8,7).Value = "02:54 PM"
$ws.Cells.Item(8,8).Value = "Early End"
$ws.Cells.Item(8,9).Value = "36 minutes early"

# Row 9: Richard Wilson
$ws.Cells.Item(9,1).Value = "Richard Wilson"
$ws.Cells.Item(9,2).Value = "TRK-1007"
$ws.Cells.Item(9,3).Value = "Downtown Construction"
$ws.Cells.Item(9,4).Value = "06:15 AM"
$ws.Cells.Item(9,5).Value = "03:15 PM"
$ws.Cells.Item(9,6).Value = "06:14 AM"
$ws.Cells.Item(9,7).Value = "03:18 PM"
$ws.Cells.Item(9,8).Value = "On Time"
$ws.Cells.Item(9,9).Value = ""

# Row 10: Joseph Moore
$ws.Cells.Item(10,1).Value = "Joseph Moore"
$ws.Cells.Item(10,2).Value = "TRK-1008"
$ws.Cells.Item(10,3).Value = "Richardson Development"
$ws.Cells.Item(10,4).Value = "06:30 AM"
$ws.Cells.Item(10,5).Value = "03:30 PM"
$ws.Cells.Item(10,6).Value = "06:55 AM"
$ws.Cells.Item(10,7).Value = "03:25 PM"
$ws.Cells.Item(10,8).Value = "Late"
$ws.Cells.Item(10,9).Value = "25 minutes late"

# Row 11: Thomas Taylor
$ws.Cells.Item(11,1).Value = "Thomas Taylor"
$ws.Cells.Item(11,2).Value = "TRK-1009"
$ws.Cells.Item(11,3).Value = "North Dallas Site"
$ws.Cells.Item(11,4).Value = "07:00 AM"
$ws.Cells.Item(11,5).Value = "04:00 PM"
$ws.Cells.Item(11,6).Value = "06:43 AM"
$ws.Cells.Item(11,7).Value = "03:52 PM"
$ws.Cells.Item(11,8).Value = "On Time"
$ws.Cells.Item(11,9).Value = ""

# Row 12: John Smith
$ws.Cells.Item(12,1).Value = "John Smith"
$ws.Cells.Item(12,2).Value = "TRK-1001"
$ws.Cells.Item(12,3).Value = "North Dallas Site"
$ws.Cells.Item(12,4).Value = "06:30 AM"
$ws.Cells.Item(12,5).Value = "03:30 PM"
$ws.Cells.Item(12,6).Value = "06:25 AM"
$ws.Cells.Item(12,7).Value = "03:35 PM"
$ws.Cells.Item(12,8).Value = "On Time"
$ws.Cells.Item(12,9).Value = ""

# Row 13: David Miller
$ws.Cells.Item(13,1).Value = "David Miller"
$ws.Cells.Item(13,2).Value = "TRK-1006"
$ws.Cells.Item(13,3).Value = "West Plano Project"
$ws.Cells.Item(13,4).Value = "06:45 AM"
$ws.Cells.Item(13,5).Value = "03:45 PM"
$ws.Cells.Item(13,6).Value = "06:47 AM"
$ws.Cells.Item(13,7).Value = "03:42 PM"
$ws.Cells.Item(13,8).Value = "On Time"
$ws.Cells.Item(13,9).Value = ""

# Row 14: Matthew Harris
$ws.Cells.Item(14,1).Value = "Matthew Harris"
$ws.Cells.Item(14,2).Value = "TRK-1013"
$ws.Cells.Item(14,3).Value = "North Dallas Site"
$ws.Cells.Item(14,4).Value = "06:45 AM"
$ws.Cells.Item(14,5).Value = "03:45 PM"
$ws.Cells.Item(14,6).Value = "06:52 AM"
$ws.Cells.Item(14,7).Value = "03:37 PM"
$ws.Cells.Item(14,8).Value = "On Time"
$ws.Cells.Item(14,9).Value = ""

# Row 15: Michael Johnson
$ws.Cells.Item(15,1).Value = "Michael Johnson"
$ws.Cells.Item(15,2).Value = "TRK-1002"
$ws.Cells.Item(15,3).Value = "West Plano Project"
$ws.Cells.Item(15,4).Value = "06:00 AM"
$ws.Cells.Item(15,5).Value = "03:00 PM"
$ws.Cells.Item(15,6).Value = "06:18 AM"
$ws.Cells.Item(15,7).Value = "03:07 PM"
$ws.Cells.Item(15,8).Value = "Late"
$ws.Cells.Item(15,9).Value = "18 minutes late"

# Row 16: Robert Williams
$ws.Cells.Item(16,1).Value = "Robert Williams"
$ws.Cells.Item(16,2).Value = "TRK-1003"
$ws.Cells.Item(16,3).Value = "Downtown Construction"
$ws.Cells.Item(16,4).Value = "07:00 AM"
$ws.Cells.Item(16,5).Value = "04:00 PM"
$ws.Cells.Item(16,6).Value = "07:28 AM"
$ws.Cells.Item(16,7).Value = "03:45 PM"
$ws.Cells.Item(16,8).Value = "Late"
$ws.Cells.Item(16,9).Value = "28 minutes late"
